$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns used for numeric-looking price/percentage strings
# remain stored as text (matching original inlineStr cells), not auto-converted
# to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.052.28"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "1.901.13"
$ws.Range("E3").Value = "  +1.80%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "312.71"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.5077"
$ws.Range("E7").Value = "  +0.66%  "
$ws.Range("D8").Value = "0.3925"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "0.09245"
$ws.Range("E9").Value = "  -4.06%  "
$ws.Range("D10").Value = "1.135"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  +2.23%  "
$ws.Range("D12").Value = "6.364"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").Value = "20.77"
$ws.Range("E13").Value = "  -0.93%  "
$ws.Range("D14").Value = "1.897.10"
$ws.Range("E14").Value = "  +1.00%  "
$ws.Range("D15").Value = "1.001"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "7.294"
$ws.Range("E16").Value = "  -1.61%  "
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").Value = "92.34"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "0.06584"
$ws.Range("E19").Value = "  -0.63%  "
$ws.Range("D20").Value = "17.77"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "6.217"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").Value = "28.110.44"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("D24").Value = "11.35"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "2.314"
$ws.Range("E25").Value = "  +1.40%  "
$ws.Range("D26").Value = "2.608"
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("D27").Value = "2.116.24"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").Value = "20.92"
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").Value = "157.37"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "127.14"
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "1.086"
$ws.Range("E31").Value = "  +1.80%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "0.1068"
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").Value = "5.607"
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "9.600"
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "0.06644"
$ws.Range("E36").Value = "  -1.52%  "
$ws.Range("D37").Value = "0.02407"
$ws.Range("E37").Value = "  +0.99%  "
$ws.Range("D38").Value = "0.2169"
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("D39").Value = "1.223"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("D40").Value = "1.262"
$ws.Range("E40").Value = "  +6.84%  "
$ws.Range("D41").Value = "0.6369"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("D42").Value = "4.989"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "11.40"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "1.001"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "13.29"
$ws.Range("E45").Value = "  -2.32%  "
$ws.Range("D46").Value = "0.5982"
$ws.Range("E46").Value = "  -0.52%  "
$ws.Range("D47").Value = "3.705"
$ws.Range("E47").Value = "  +1.18%  "
$ws.Range("D48").Value = "1.274"
$ws.Range("E48").Value = "  +0.77%  "
$ws.Range("D49").Value = "2.011"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").Value = "122.56"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "1.178"
$ws.Range("E51").Value = "  -1.34%  "

Write-Host "Applied cryptos update"
